# [Add] - doc chuc nang phan mem
#
# The title block of the document starts with two centered, bold, 44pt
# paragraphs:
#   1) "test"                                   (+ the hidden _GoBack bookmark)
#   2) "NHAP MON CONG NGHE PHAN MEM"
#
# The edit removes the stray "test" paragraph by deleting its text and
# merging it into the following paragraph, so the _GoBack bookmark (which
# sits between the two) ends up at the start of the remaining title
# paragraph, immediately before the real title run.

$d = $word.ActiveDocument

# Locate the standalone "test" run unambiguously.
$rng = $d.Content
$found = $rng.Find.Execute("test", $true, $true, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $hostPara = $rng.Paragraphs(1)

    # Merge the "test" paragraph into the next one first, by deleting its
    # trailing paragraph mark. Doing the merge *before* clearing the text
    # keeps the (hidden) _GoBack bookmark that sits right after "test"
    # anchored correctly at the join point instead of being swallowed.
    $paraEnd = $hostPara.Range.End
    $mark = $d.Range($paraEnd - 1, $paraEnd)
    $mark.Delete()

    # Now remove the "test" text itself; the run collapses away with it.
    $rng.Delete()
}
